$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.927.61'
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = '1.879.80'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("E4").Value = '  +1.48%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.23'
$ws.Range("E5").Value = '  +0.75%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4679'
$ws.Range("E7").Value = '  -1.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3912'
$ws.Range("E8").Value = '  -1.22%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.74'
$ws.Range("E9").Value = '  -2.20%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07941'
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.007'
$ws.Range("E11").Value = '  -1.51%  '
$ws.Range("E12").Value = '  -1.38%  '
$ws.Range("D13").Value = '1.910.92'
$ws.Range("E13").Value = '  +1.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.943'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.098'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.06789'
$ws.Range("E17").Value = '  +2.33%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '87.42'
$ws.Range("E18").Value = '  +0.23%  '
$ws.Range("E19").Value = '  -0.56%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.95'
$ws.Range("E20").Value = '  -1.48%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.018'
$ws.Range("E21").Value = '  +1.38%  '
$ws.Range("D22").Value = '27.944.27'
$ws.Range("E22").Value = '  +0.81%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.452'
$ws.Range("E23").Value = '  -0.78%  '
$ws.Range("E24").Value = '  -0.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.357'
$ws.Range("E25").Value = '  +2.47%  '
$ws.Range("D26").Value = '2.132.98'
$ws.Range("E26").Value = '  +1.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.27'
$ws.Range("E27").Value = '  +1.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.93'
$ws.Range("E28").Value = '  -1.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.066'
$ws.Range("E29").Value = '  -1.80%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.448'
$ws.Range("E30").Value = '  -2.28%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '120.57'
$ws.Range("E31").Value = '  -1.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09524'
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9543'
$ws.Range("E33").Value = '  -1.59%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.661'
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.317'
$ws.Range("E35").Value = '  +0.41%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.345'
$ws.Range("E36").Value = '  -7.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06111'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02235'
$ws.Range("E38").Value = '  -1.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.203'
$ws.Range("E39").Value = '  -2.08%  '
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.094'
$ws.Range("E41").Value = '  -1.05%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5867'
$ws.Range("E42").Value = '  -2.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1891'
$ws.Range("E43").Value = '  -1.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.15'
$ws.Range("E44").Value = '  -1.15%  '
$ws.Range("E45").Value = '  +1.65%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5622'
$ws.Range("E46").Value = '  -1.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.16'
$ws.Range("E47").Value = '  -1.10%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.403'
$ws.Range("E48").Value = '  -0.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.914'
$ws.Range("E49").Value = '  -1.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06858'
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '113.38'
$ws.Range("E51").Value = '  +0.46%  '
